# DOMA-3097 review fixes for contact export
#
# The "Unit Type" column used to be the last column (F) of the contact
# import template. Review feedback moved it to be the 3rd column (C),
# right after "Unit Name", pushing "Phones", "Full Name" and "Email" one
# column to the right (C->D, D->E, E->F). The mailto: hyperlink that lived
# on the "Email" cell has to follow that cell from E2 to F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move column F ("Unit Type") so it becomes column C; C,D,E shift right to
# D,E,F respectively. This single Cut + Insert reproduces the header
# re-order, the row-3..50 placeholder-cell style shuffle, and the
# column-width shuffle all in one shot (Excel carries cell content AND
# formatting along with a cut column).
$ws.Columns("F").Cut()
$ws.Columns("C").Insert()

# The cut/insert above does not drag the worksheet's hyperlink along with
# the "Email" cell, so re-point it by hand: drop the stale one (still
# anchored on the old E2) and recreate it on the new location, F2.
$ws.Range("E2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:test@example.com", "", "", "test@example.com")

# Hyperlinks.Add stamps the target cell with the builtin blue/underlined
# "Hyperlink" style; the source file does not use that style here, so
# restore F2's original (non-hyperlink) formatting by pasting just the
# formats back from another cell that already carries that plain style.
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
